$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BM")

$ws.Range("B3").Value = 0.2302
$ws.Range("C3").Value = 0.229
$ws.Range("D3").Value = 0.0668
$ws.Range("E3").Value = 0.0397
$ws.Range("F3").Value = 0.2292
$ws.Range("G3").Value = 0.0389
$ws.Range("H3").Value = 0.2302
$ws.Range("B4").Value = 0.229
$ws.Range("C4").Value = 0.2286
$ws.Range("D4").Value = 0.0482
$ws.Range("E4").Value = 0.0491
$ws.Range("F4").Value = 0.2296
$ws.Range("G4").Value = 0.0428
$ws.Range("H4").Value = 0.2296
$ws.Range("B5").Value = 0.2305
$ws.Range("C5").Value = 0.2284
$ws.Range("D5").Value = 0.0393
$ws.Range("E5").Value = 0.0471
$ws.Range("F5").Value = 0.2306
$ws.Range("G5").Value = 0.0482
$ws.Range("H5").Value = 0.2298
$ws.Range("B6").Value = 0.2308
$ws.Range("C6").Value = 0.2313
$ws.Range("D6").Value = 0.04
$ws.Range("E6").Value = 0.0567
$ws.Range("F6").Value = 0.2313
$ws.Range("G6").Value = 0.0463
$ws.Range("H6").Value = 0.2322
$ws.Range("B7").Value = 0.2316
$ws.Range("C7").Value = 0.2312
$ws.Range("D7").Value = 0.0441
$ws.Range("E7").Value = 0.0439
$ws.Range("F7").Value = 0.2318
$ws.Range("G7").Value = 0.048
$ws.Range("H7").Value = 0.2326
$ws.Range("B8").Value = 0.2328
$ws.Range("C8").Value = 0.2327
$ws.Range("D8").Value = 0.0479
$ws.Range("E8").Value = 0.0485
$ws.Range("F8").Value = 0.2322
$ws.Range("G8").Value = 0.0476
$ws.Range("H8").Value = 0.2326
$ws.Range("B9").Value = 0.233
$ws.Range("C9").Value = 0.2328
$ws.Range("D9").Value = 0.0499
$ws.Range("E9").Value = 0.0518
$ws.Range("F9").Value = 0.2327
$ws.Range("G9").Value = 0.0525
$ws.Range("H9").Value = 0.2336
$ws.Range("B10").Value = 0.2333
$ws.Range("C10").Value = 0.2334
$ws.Range("D10").Value = 0.059
$ws.Range("E10").Value = 0.0614
$ws.Range("F10").Value = 0.2334
$ws.Range("G10").Value = 0.0584
$ws.Range("H10").Value = 0.2334
$ws.Range("B11").Value = 0.2302
$ws.Range("C11").Value = 0.2307
$ws.Range("D11").Value = 0.0806
$ws.Range("E11").Value = 0.0856
$ws.Range("F11").Value = 0.2298
$ws.Range("G11").Value = 0.0788
$ws.Range("H11").Value = 0.2297
$ws.Range("B12").Value = 0.1923
$ws.Range("C12").Value = 0.1952
$ws.Range("D12").Value = 0.1438
$ws.Range("E12").Value = 0.1452
$ws.Range("F12").Value = 0.1935
$ws.Range("G12").Value = 0.1465
$ws.Range("H12").Value = 0.1923
$ws.Range("B13").Value = 0.1828
$ws.Range("C13").Value = 0.183
$ws.Range("D13").Value = 0.1548
$ws.Range("E13").Value = 0.1551
$ws.Range("F13").Value = 0.182
$ws.Range("G13").Value = 0.1542
$ws.Range("H13").Value = 0.1815
$ws.Range("B14").Value = 0.1736
$ws.Range("C14").Value = 0.1721
$ws.Range("D14").Value = 0.162
$ws.Range("E14").Value = 0.1645
$ws.Range("F14").Value = 0.1756
$ws.Range("G14").Value = 0.1641
$ws.Range("H14").Value = 0.1707
$ws.Range("B15").Value = 0.1622
$ws.Range("C15").Value = 0.1639
$ws.Range("D15").Value = 0.1727
$ws.Range("E15").Value = 0.1737
$ws.Range("F15").Value = 0.1631
$ws.Range("G15").Value = 0.1748
$ws.Range("H15").Value = 0.1626
$ws.Range("B16").Value = 0.16
$ws.Range("C16").Value = 0.1604
$ws.Range("D16").Value = 0.1825
$ws.Range("E16").Value = 0.1827
$ws.Range("F16").Value = 0.16
$ws.Range("G16").Value = 0.184
$ws.Range("H16").Value = 0.161
$ws.Range("B17").Value = 0.1667
$ws.Range("C17").Value = 0.165
$ws.Range("D17").Value = 0.1916
$ws.Range("E17").Value = 0.193
$ws.Range("F17").Value = 0.1635
$ws.Range("G17").Value = 0.1924
$ws.Range("H17").Value = 0.1678
$ws.Range("B18").Value = 0.1793
$ws.Range("C18").Value = 0.1807
$ws.Range("D18").Value = 0.2007
$ws.Range("E18").Value = 0.2008
$ws.Range("F18").Value = 0.18
$ws.Range("G18").Value = 0.2007
$ws.Range("H18").Value = 0.1824
$ws.Range("B19").Value = 0.1963
$ws.Range("C19").Value = 0.1967
$ws.Range("D19").Value = 0.2078
$ws.Range("E19").Value = 0.2078
$ws.Range("F19").Value = 0.1968
$ws.Range("G19").Value = 0.2079
$ws.Range("H19").Value = 0.1993
$ws.Range("B20").Value = 0.2103
$ws.Range("C20").Value = 0.2104
$ws.Range("D20").Value = 0.2127
$ws.Range("E20").Value = 0.2126
$ws.Range("F20").Value = 0.2097
$ws.Range("G20").Value = 0.2127
$ws.Range("H20").Value = 0.2106
$ws.Range("B21").Value = 0.2161
$ws.Range("C21").Value = 0.2161
$ws.Range("D21").Value = 0.2156
$ws.Range("E21").Value = 0.2156
$ws.Range("F21").Value = 0.216
$ws.Range("G21").Value = 0.2156
$ws.Range("H21").Value = 0.2161

$excel.Calculate()

$ws.Range("H21").Select()
